# Append new ticker rows to the end of the sheet (rows 414-418)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("TAO-USD", "IMX-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 414
for ($i = 0; $i -lt $newTickers.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
